# "Work: tests for potenot calculator"
# - Journal sheet: log a new work entry (row 8) for PotenotTaskServece tests,
#   and move the view/selection.
# - "PotenotTask test data" sheet: fill in newly-computed columns (V/W/X/Y,
#   R/S) for the potenot-calculator test rows, tweak a couple of existing
#   results, and add two new "actual result" columns (AE/AF) while clearing
#   the now-superseded AA values.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Journal")
$ws2 = $wb.Worksheets.Item("PotenotTask test data")

# ---------------------------------------------------------------------
# Journal: new row describing the work done (tests for PotenotService)
# ---------------------------------------------------------------------
$ws1.Range("A8").Value = "Разработка тестов для PotenotService"
$ws1.Range("B8").Value = 45547
$ws1.Range("C8").Value = 0.0590277777777778
$ws1.Range("D8").Value = "PotenotTaskServece"

# ---------------------------------------------------------------------
# PotenotTask test data: row 3 (target2) -- fill second/third angle pair
# ---------------------------------------------------------------------
$ws2.Range("V3").Value = 197.2045
$ws2.Range("W3").Value = 3.444334557
$ws2.Range("X3").Value = 242.3223
$ws2.Range("Y3").Value = 4.232825831
$ws2.Range("Z3").Value = 18514

# ---------------------------------------------------------------------
# PotenotTask test data: row 4 (target3)
# ---------------------------------------------------------------------
$ws2.Range("R4").Value = -630177.7442
$ws2.Range("S4").Value = 683001.4427
$ws2.Range("T4").Value = 317.1813
$ws2.Range("U4").Value = 5.537992742
$ws2.Range("V4").Value = 324.3454
$ws2.Range("W4").Value = 5.665018775
$ws2.Range("X4").Value = 325.5748
$ws2.Range("Y4").Value = 5.689133407
$ws2.Range("Z4").Value = -709005
$ws2.Range("AA4").ClearContents()
$ws2.Range("AE4").Value = -709005
$ws2.Range("AF4").Value = 846434

# ---------------------------------------------------------------------
# PotenotTask test data: row 5 (target4)
# ---------------------------------------------------------------------
$ws2.Range("N5").Value = 0.363236954
$ws2.Range("R5").Value = 49548.54946
$ws2.Range("S5").Value = -163086.6288
$ws2.Range("T5").Value = 163.0601
$ws2.Range("U5").Value = 2.846636858
$ws2.Range("V5").Value = 183.5444
$ws2.Range("W5").Value = 3.209873812
$ws2.Range("X5").Value = 206.2047
$ws2.Range("Y5").Value = 3.601423886
$ws2.Range("Z5").Value = 183368
$ws2.Range("AA5").ClearContents()
$ws2.Range("AE5").Value = 15295
$ws2.Range("AF5").Value = -176044

# ---------------------------------------------------------------------
# View state: scroll/selection on both sheets. "PotenotTask test data" is
# touched first so the final Activate() below leaves "Journal" as the
# selected tab, matching the original tab order.
# ---------------------------------------------------------------------
$ws2.Activate()
$win2 = $excel.ActiveWindow
$win2.ScrollColumn = 21
$win2.ScrollRow = 1
$ws2.Range("AC17").Select()

$ws1.Activate()
$win1 = $excel.ActiveWindow
$win1.ScrollColumn = 2
$win1.ScrollRow = 1
$ws1.Range("D18").Select()

# Best-effort: default column width tweaks from the diff (engine may not
# persist StandardWidth on export, but set it for correctness anyway).
$ws1.StandardWidth = 11.6875
$ws2.StandardWidth = 8.75
